$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: title-case standalone Spanish articles/prepositions within a name ---
function Fix-SpanishName($s) {
    if ($s -eq $null) { return $s }
    $preps = @("de", "del", "el", "la", "las", "lo", "los", "y", "en", "un", "una")
    $parts = $s -split " "
    $outParts = @()
    foreach ($p in $parts) {
        $isPrep = $false
        foreach ($pr in $preps) {
            if ($p -ceq $pr) {
                $isPrep = $true
            }
        }
        if ($isPrep) {
            $outParts += ($p.Substring(0, 1).ToUpper() + $p.Substring(1))
        } else {
            $outParts += $p
        }
    }
    return [string]::Join(" ", $outParts)
}

# --- 1. Remove the trailing footnote rows (971:975) ---
$ws.Rows("971:975").Delete()

# --- 2. Re-title-case the state (col A) / municipality (col B) names in the data rows ---
$lastRow = 969
for ($r = 2; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $valA = $cellA.Value2
    if ($valA -ne $null -and $valA -ne "") {
        $fixedA = Fix-SpanishName($valA)
        $cellA.Value = $fixedA
    }

    $cellB = $ws.Cells.Item($r, 2)
    $valB = $cellB.Value2
    if ($valB -ne $null -and $valB -ne "") {
        $fixedB = Fix-SpanishName($valB)
        $cellB.Value = $fixedB
    }
}

# --- 3. Rename the header row to the snake_case column names ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"
